$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"

$ws.Range("C4").Select()
